# Update test case data on sheet "three" (3rd worksheet), rows 2-6, columns A-C.
# New shared strings must be introduced in this exact order so they are appended
# to the shared string table as: Catheryn, Spring, Troy, Camilla, Trent, Perry,
# Norbert, Mickey, Jerrod, Ron, Gidget, Lynna, Ena, Jeannetta.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$ws.Range("A4").Value = "Catheryn"
$ws.Range("C4").Value = "Spring"
$ws.Range("A5").Value = "Troy"
$ws.Range("B5").Value = "Camilla"
$ws.Range("C5").Value = "Trent"
$ws.Range("A3").Value = "Perry"
$ws.Range("B3").Value = "Norbert"
$ws.Range("C3").Value = "Mickey"
$ws.Range("A6").Value = "Jerrod"
$ws.Range("B6").Value = "Ron"
$ws.Range("C6").Value = "Gidget"
$ws.Range("A2").Value = "Lynna"
$ws.Range("B2").Value = "Ena"
$ws.Range("C2").Value = "Jeannetta"

$ws.Range("B4").Value = "Julio"
